# WI - Cara Packing Board Menggunakan Tali Rafia.docx
#
# 1) Center the vertical alignment of the 5 "value" cells (4316 dxa wide,
#    second column) in the header table.
# 2) Split the document-number run "01/No.13 - Dok.03/2022" so that a new
#    "VST/" segment is inserted before "2022", and move the (single)
#    "_GoBack" bookmark from its old home (in front of "Jika, Tidak
#    Memungkinkan untuk ditali :") to sit right after the new "VST/" text.

$d = $word.ActiveDocument

# --- 1) vAlign = center on the 5 header-table value cells -----------------
$t = $d.Tables.Item(1)
for ($i = 1; $i -le 5; $i++) {
    $cell = $t.Cell($i, 2)
    $cell.VerticalAlignment = 1   # wdCellAlignVerticalCenter
}

# --- 2) Split "01/No.13 - Dok.03/2022" into 3 runs + relocate _GoBack -----

# Locate "01/No.13" (start of the document-number value) ...
$findDocNum = $d.Content
$foundDocNum = $findDocNum.Find.Execute("01/No.13", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pDocNum = $findDocNum.Start

# ... and "2022" - this is where "VST/" needs to be inserted in front of.
$find2022 = $d.Content
$found2022 = $find2022.Find.Execute("2022", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p2022 = $find2022.Start

# Insert the new "VST/" text right before "2022".
$insertRng = $d.Range($p2022, $p2022)
$insertRng.InsertBefore("VST/")

# Re-establish the run break between "No Dokumen : " and "01/No.13 - Dok.03/"
# (a plain insert coalesces same-formatted neighbouring runs together, so
# nudge the formatting off and back on to force the boundary to stick).
$docNumRng = $d.Range($pDocNum, $p2022)
$docNumRng.Bold = 1
$docNumRng.Bold = 0

# Likewise force "VST/" to live in its own run rather than merge with the
# text before or after it.
$vstRng = $d.Range($p2022, $p2022 + 4)
$vstRng.Bold = 1
$vstRng.Bold = 0

# Move the "_GoBack" bookmark so it now sits between "VST/" and "2022".
# Adding a bookmark with a name that already exists elsewhere relocates it,
# removing the old bookmarkStart/bookmarkEnd pair automatically.
$bmRng = $d.Range($p2022 + 4, $p2022 + 4)
$bmRng.Bookmarks.Add("_GoBack")
